$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gains two new columns ("Hydro interest rate" / "Hydro lifetime
# (years)") inserted right after the Wind columns (between the old columns H
# and I). Copy the existing Plant/Infrastructure block (H1:K2) two columns to
# the right (-> J1:M2), then clear the vacated H1:I2 cells so they can hold
# the new Hydro headers. Using copy/paste (rather than a full column insert)
# keeps the pre-existing column width formatting anchored to its original
# column position, exactly like the source edit.
$ws.Range("H1:K2").Copy()
$ws.Range("J1").PasteSpecial()
$ws.Range("H1:I2").ClearContents()
$excel.CutCopyMode = $false

# New header cells for the inserted Hydro columns (row 2 stays blank for
# this country, same as the source edit).
$ws.Range("H1").Value = "Hydro interest rate"
$ws.Range("I1").Value = "Hydro lifetime (years)"

# Column I (the "Hydro lifetime (years)" header) is widened/best-fit to the
# new text; auto-fit it to the header, nudging to the closest width the
# host's column-width quantization can represent.
$ws.Range("I1").EntireColumn.AutoFit()
$ws.Range("I1").EntireColumn.ColumnWidth = 18.3

# Update the active selection to rest on I2, matching the source edit.
$ws.Range("I2").Select()
